$wb = $excel.ActiveWorkbook

# Sheet ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2929.625
$ws.Range("I64").Value = 3130.3333
$ws.Range("K64").Value = 3130.3333
$ws.Range("M64").Value = -2882.3333

# Sheet ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2929.625
$ws.Range("I67").Value = 3130.3333
$ws.Range("K67").Value = 3130.3333
$ws.Range("M67").Value = -2272.3333

# Sheet ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 5500
$ws.Range("I70").Value = 4833.3335
$ws.Range("J70").Value = 9500
$ws.Range("K70").Value = 14500.0005
$ws.Range("L70").Value = 28500
$ws.Range("M70").Value = -14230.0005
$ws.Range("N70").Value = -29040

# Sheet ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 5500
$ws.Range("I73").Value = 4833.3335
$ws.Range("J73").Value = 9500
$ws.Range("K73").Value = 14500.0005
$ws.Range("L73").Value = 28500
$ws.Range("M73").Value = -13564.0005
$ws.Range("N73").Value = -30372

# Sheet ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 29999.5
$ws.Range("I106").Value = 29999.5
$ws.Range("K106").Value = 29999.5
$ws.Range("M106").Value = -29368.5

# Sheet ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3124.75
$ws.Range("I116").Value = 2666.3333
$ws.Range("J116").Value = 4500
$ws.Range("K116").Value = 2666.3333
$ws.Range("L116").Value = 4500
$ws.Range("M116").Value = 775.6667000000002
$ws.Range("N116").Value = -11384

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2179.0908
$ws.Range("I137").Value = 1446.8
$ws.Range("K137").Value = 4340.4
$ws.Range("M137").Value = -1790.4

# Sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2009.1765
$ws.Range("I138").Value = 1537.4
$ws.Range("J138").Value = 2381.6316
$ws.Range("K138").Value = 4612.200000000001
$ws.Range("L138").Value = 7144.8948
$ws.Range("M138").Value = 527.7999999999993
$ws.Range("N138").Value = -17424.8948

# Sheet ARM row 24
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").ClearContents()
$ws.Range("N24").Value = 0

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5932.6665
$ws.Range("I32").Value = 5282.273
$ws.Range("K32").Value = 5282.273
$ws.Range("M32").Value = -4995.273

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2200.5217
$ws.Range("I45").Value = 2025.75
$ws.Range("K45").Value = 2025.75
$ws.Range("M45").Value = -1648.75

# Sheet ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = 0

# Sheet ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = 0

# Sheet ARM row 100
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

# Sheet ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3192.818
$ws.Range("J102").Value = 3293
$ws.Range("L102").Value = 3293
$ws.Range("N102").Value = -6537

# Sheet ARM row 118
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 80408
$ws.Range("J118").Value = 80408
$ws.Range("L118").Value = 80408
$ws.Range("N118").Value = -83722

# Sheet ARM row 119
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 9999.125
$ws.Range("J119").Value = 9999.125
$ws.Range("L119").Value = 9999.125
$ws.Range("N119").Value = -19675.125

# Sheet ARM row 130
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 51610.5
$ws.Range("J130").Value = 51610.5
$ws.Range("L130").Value = 51610.5
$ws.Range("N130").Value = -61650.5

# Sheet BSM row 6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 24440.111
$ws.Range("J6").Value = 24440.111
$ws.Range("L6").Value = 24440.111
$ws.Range("N6").Value = -24666.111

# Sheet BSM row 55
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("N55").Value = 0

# Sheet BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 7562.3
$ws.Range("I80").Value = 287.5
$ws.Range("J80").Value = 18474.5
$ws.Range("K80").Value = 287.5
$ws.Range("L80").Value = 18474.5
$ws.Range("M80").Value = 710.5
$ws.Range("N80").Value = -20470.5

# Sheet BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 7562.3
$ws.Range("I83").Value = 287.5
$ws.Range("J83").Value = 18474.5
$ws.Range("K83").Value = 1437.5
$ws.Range("L83").Value = 92372.5
$ws.Range("M83").Value = 3554.5
$ws.Range("N83").Value = -102356.5

# Sheet BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4051.5715
$ws.Range("I94").Value = 3380
$ws.Range("J94").Value = 4947
$ws.Range("K94").Value = 3380
$ws.Range("L94").Value = 4947
$ws.Range("M94").Value = -2929
$ws.Range("N94").Value = -5849

# Sheet BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3558.1428
$ws.Range("I105").Value = 3181.4
$ws.Range("K105").Value = 3181.4
$ws.Range("M105").Value = -1434.4

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1874.5
$ws.Range("I134").Value = 1874.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5623.5
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -3088.5

# Sheet CRP row 20
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 90000
$ws.Range("J20").Value = 90000
$ws.Range("L20").Value = 90000
$ws.Range("N20").Value = -90472

# Sheet CRP row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Sheet CRP row 30
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 90000
$ws.Range("J30").Value = 90000
$ws.Range("L30").Value = 90000
$ws.Range("N30").Value = -90182

# Sheet CRP row 59
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 75076.336
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

# Sheet CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1425.6666
$ws.Range("I99").Value = 1251.6364
$ws.Range("J99").Value = 1904.25
$ws.Range("K99").Value = 1251.6364
$ws.Range("L99").Value = 1904.25
$ws.Range("M99").Value = 246.3635999999999
$ws.Range("N99").Value = -4900.25

# Sheet CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1479.8334
$ws.Range("I107").Value = 1396
$ws.Range("J107").Value = 1899
$ws.Range("K107").Value = 1396
$ws.Range("L107").Value = 1899
$ws.Range("M107").Value = 524
$ws.Range("N107").Value = -5739

# Sheet CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1425.6666
$ws.Range("I126").Value = 1251.6364
$ws.Range("J126").Value = 1904.25
$ws.Range("K126").Value = 3754.9092
$ws.Range("L126").Value = 5712.75
$ws.Range("M126").Value = -1284.9092
$ws.Range("N126").Value = -10652.75

# Sheet CRP row 128
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960

# Sheet CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2700.5833
$ws.Range("I5").Value = 2491.5454
$ws.Range("K5").Value = 7474.6362
$ws.Range("M5").Value = -7362.6362

# Sheet CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 350.375
$ws.Range("J17").Value = 462
$ws.Range("L17").Value = 1386
$ws.Range("N17").Value = -1724

# Sheet CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2700.5833
$ws.Range("I135").Value = 2491.5454
$ws.Range("K135").Value = 22423.9086
$ws.Range("M135").Value = -19888.9086

# Sheet GSM row 23
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 1950

# Sheet GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11015.5
$ws.Range("I70").Value = 7431.3335
$ws.Range("K70").Value = 7431.3335
$ws.Range("M70").Value = -7161.3335

# Sheet GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11015.5
$ws.Range("I73").Value = 7431.3335
$ws.Range("K73").Value = 7431.3335
$ws.Range("M73").Value = -6495.3335

# Sheet GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4287.9287
$ws.Range("I80").Value = 3616
$ws.Range("K80").Value = 3616
$ws.Range("M80").Value = -2618

# Sheet GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4287.9287
$ws.Range("I83").Value = 3616
$ws.Range("K83").Value = 18080
$ws.Range("M83").Value = -13088

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1651.4615
$ws.Range("I102").Value = 1369
$ws.Range("K102").Value = 1369
$ws.Range("M102").Value = 253

# Sheet GSM row 104
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 33065.11
$ws.Range("J104").Value = 33065.11
$ws.Range("L104").Value = 33065.11
$ws.Range("N104").Value = -40053.11

# Sheet GSM row 114
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H114").Value = 100000
$ws.Range("J114").Value = 100000
$ws.Range("L114").Value = 100000
$ws.Range("N114").Value = -108678

# Sheet GSM row 128
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2081.5715
$ws.Range("I132").Value = 2081.5715
$ws.Range("K132").Value = 6244.7145
$ws.Range("M132").Value = -3714.7145

# Sheet LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 616.2857
$ws.Range("I16").Value = 644.9167
$ws.Range("J16").Value = 444.5
$ws.Range("K16").Value = 644.9167
$ws.Range("L16").Value = 444.5
$ws.Range("M16").Value = -474.9167
$ws.Range("N16").Value = -784.5

# Sheet LTW row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0

# Sheet LTW row 128
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 53795.8
$ws.Range("J128").Value = 53795.8
$ws.Range("L128").Value = 53795.8
$ws.Range("N128").Value = -63755.8

# Sheet WVR row 41
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17909
$ws.Range("I41").Value = 17909
$ws.Range("K41").Value = 17909
$ws.Range("M41").Value = -17519

# Sheet WVR row 124
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 23970.166
$ws.Range("J124").Value = 23970.166
$ws.Range("L124").Value = 23970.166
$ws.Range("N124").Value = -33790.166
